$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New student records to append (case study cuoi khoa)
$newRows = @(
    @("SV017", "Nguyễn Quốc Hưng", 1998, "Nam", "Tốt nghiệp"),
    @("SV018", "Nguyễn Nghĩa", 1996, "Nam", "Tốt nghiệp"),
    @("SV019", "Trần Hồng", 1996, "Nữ", "Chưa tốt nghiệp")
)

$startRow = 18
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
